$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 50009.332
$ws.Range("I21").Value = 58011.4
$ws.Range("J21").Value = 9999
$ws.Range("K21").Value = 58011.4
$ws.Range("L21").Value = 9999
$ws.Range("M21").Value = -57543.4
$ws.Range("N21").Value = -10935

$ws.Range("H23").Value = 50009.332
$ws.Range("I23").Value = 58011.4
$ws.Range("J23").Value = 9999
$ws.Range("K23").Value = 58011.4
$ws.Range("L23").Value = 9999
$ws.Range("M23").Value = -57777.4
$ws.Range("N23").Value = -10467

$ws.Range("H137").Value = 1927.25
$ws.Range("I137").Value = 1298.1333
$ws.Range("J137").Value = 2376.6191
$ws.Range("K137").Value = 3894.3999
$ws.Range("L137").Value = 7129.8573
$ws.Range("M137").Value = -1344.3999
$ws.Range("N137").Value = -12229.8573

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3236.7817
$ws.Range("I61").Value = 3954.0732
$ws.Range("K61").Value = 3954.0732
$ws.Range("M61").Value = -3742.0732

$ws.Range("H74").Value = 1796.909
$ws.Range("I74").Value = 1502.3334
$ws.Range("J74").Value = 2150.4
$ws.Range("K74").Value = 1502.3334
$ws.Range("L74").Value = 2150.4
$ws.Range("M74").Value = -628.3334
$ws.Range("N74").Value = -3898.4

$ws.Range("H77").Value = 1796.909
$ws.Range("I77").Value = 1502.3334
$ws.Range("J77").Value = 2150.4
$ws.Range("K77").Value = 7511.666999999999
$ws.Range("L77").Value = 10752
$ws.Range("M77").Value = -3143.666999999999
$ws.Range("N77").Value = -19488

$ws.Range("H136").Value = 3236.7817
$ws.Range("I136").Value = 3954.0732
$ws.Range("K136").Value = 11862.2196
$ws.Range("M136").Value = -9312.2196

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("N46").ClearContents()

$ws.Range("H134").Value = 7040.4585
$ws.Range("I134").Value = 11834.454
$ws.Range("J134").Value = 2984
$ws.Range("K134").Value = 35503.362
$ws.Range("L134").Value = 8952
$ws.Range("M134").Value = -32968.362
$ws.Range("N134").Value = -14022

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3609.0366
$ws.Range("I31").Value = 1830.5
$ws.Range("J31").Value = 4261.1665
$ws.Range("K31").Value = 1830.5
$ws.Range("L31").Value = 4261.1665
$ws.Range("M31").Value = -1535.5
$ws.Range("N31").Value = -4851.1665

$ws.Range("H34").Value = 3609.0366
$ws.Range("I34").Value = 1830.5
$ws.Range("J34").Value = 4261.1665
$ws.Range("K34").Value = 1830.5
$ws.Range("L34").Value = 4261.1665
$ws.Range("M34").Value = -1628.5
$ws.Range("N34").Value = -4665.1665

$ws.Range("H58").Value = 1629.6
$ws.Range("I58").Value = 954.44446
$ws.Range("J58").Value = 2642.3333
$ws.Range("K58").Value = 954.44446
$ws.Range("L58").Value = 2642.3333
$ws.Range("M58").Value = -751.44446
$ws.Range("N58").Value = -3048.3333

$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H132").Value = 5804.364
$ws.Range("I132").Value = 10000
$ws.Range("J132").Value = 5384.8
$ws.Range("K132").Value = 30000
$ws.Range("L132").Value = 16154.4
$ws.Range("M132").Value = -27470
$ws.Range("N132").Value = -21214.4

$ws.Range("H134").Value = 2414.2778
$ws.Range("I134").Value = 2539.394
$ws.Range("J134").Value = 1038
$ws.Range("K134").Value = 7618.181999999999
$ws.Range("L134").Value = 3114
$ws.Range("M134").Value = -5083.181999999999
$ws.Range("N134").Value = -8184

$ws.Range("H136").Value = 1629.6
$ws.Range("I136").Value = 954.44446
$ws.Range("J136").Value = 2642.3333
$ws.Range("K136").Value = 2863.33338
$ws.Range("L136").Value = 7926.999899999999
$ws.Range("M136").Value = -313.33338
$ws.Range("N136").Value = -13026.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2732.4614
$ws.Range("I68").Value = 4188.7744
$ws.Range("J68").Value = 1771.9149
$ws.Range("K68").Value = 12566.3232
$ws.Range("L68").Value = 5315.7447
$ws.Range("M68").Value = -11755.3232
$ws.Range("N68").Value = -6937.7447

$ws.Range("H71").Value = 2732.4614
$ws.Range("I71").Value = 4188.7744
$ws.Range("J71").Value = 1771.9149
$ws.Range("K71").Value = 37698.9696
$ws.Range("L71").Value = 15947.2341
$ws.Range("M71").Value = -33642.9696
$ws.Range("N71").Value = -24059.2341

$ws.Range("H80").Value = 6250
$ws.Range("I80").Value = 2500
$ws.Range("J80").Value = 7000
$ws.Range("K80").Value = 7500
$ws.Range("L80").Value = 21000
$ws.Range("M80").Value = -6564
$ws.Range("N80").Value = -22872

$ws.Range("H83").Value = 6250
$ws.Range("I83").Value = 2500
$ws.Range("J83").Value = 7000
$ws.Range("K83").Value = 22500
$ws.Range("L83").Value = 63000
$ws.Range("M83").Value = -17820
$ws.Range("N83").Value = -72360

$ws.Range("H131").Value = 14103737
$ws.Range("I131").Value = 6667183.5
$ws.Range("J131").Value = 15874346
$ws.Range("K131").Value = 20001550.5
$ws.Range("L131").Value = 47623038
$ws.Range("M131").Value = -19996510.5
$ws.Range("N131").Value = -47633118

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 30964616
$ws.Range("I132").Value = 43348464
$ws.Range("K132").Value = 130045392
$ws.Range("M132").Value = -130042862

$ws.Range("H136").Value = 5464.7817
$ws.Range("J136").Value = 8193.267
$ws.Range("L136").Value = 24579.801
$ws.Range("N136").Value = -29679.801

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

$ws.Range("H125").Value = 42429
$ws.Range("J125").Value = 42429
$ws.Range("L125").Value = 42429
$ws.Range("N125").Value = -52269

$ws.Range("H132").Value = 1941.6765
$ws.Range("I132").Value = 1786.1305
$ws.Range("J132").Value = 2266.9092
$ws.Range("K132").Value = 5358.3915
$ws.Range("L132").Value = 6800.7276
$ws.Range("M132").Value = -2828.3915
$ws.Range("N132").Value = -11860.7276

$ws.Range("H136").Value = 3198.0344
$ws.Range("I136").Value = 3392.9
$ws.Range("J136").Value = 2765
$ws.Range("K136").Value = 10178.7
$ws.Range("L136").Value = 8295
$ws.Range("M136").Value = -7628.700000000001
$ws.Range("N136").Value = -13395

$ws.Range("H139").Value = 61999.75
$ws.Range("J139").Value = 61999.75
$ws.Range("L139").Value = 61999.75
$ws.Range("N139").Value = -72279.75
